$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: new date entry with Additional Effort only
$ws.Range("A10").Value = 41177
$ws.Range("A10").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("C10").Value = 1

# Row 11: new date entry with Effort and a comment
$ws.Range("A11").Value = 41178
$ws.Range("A11").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B11").Value = 2.5
$ws.Range("D11").Value = "waitForEvent, setEvent implemented. TC03 added, but does not yet run"

# Update the selection to match the new last-used cell
$ws.Range("D11").Select()
